# Expand rows for chart 22 ("Access to appropriate information and advice")
# The existing row (n=22) combined both survey questions q154 and q155 into a
# single "horacio_id" value ("q154,q155") with no explicit id/topic. This
# change splits that single row into two rows - one per question - each
# carrying its own horacio_id and id, sharing the same topic.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 currently holds the combined "q154,q155" entry (n=22). Insert a
# fresh row right after it to hold the second question's data; the sheet's
# used range grows from M143 to M144 and every following row shifts down
# by one.
$ws.Rows.Item(24).Insert()

# Copy the shared (chapter/section/subsection/n/question) values from row 23
# into the new row 24 so both rows describe the same chart/topic.
$ws.Cells.Item(24, 1).Value = $ws.Cells.Item(23, 1).Value2
$ws.Cells.Item(24, 2).Value = $ws.Cells.Item(23, 2).Value2
$ws.Cells.Item(24, 3).Value = $ws.Cells.Item(23, 3).Value2
$ws.Cells.Item(24, 4).Value = $ws.Cells.Item(23, 4).Value2
$ws.Cells.Item(24, 5).Value = $ws.Cells.Item(23, 5).Value2
$ws.Cells.Item(24, 9).Value = $ws.Cells.Item(23, 9).Value2
$ws.Cells.Item(24, 10).Value = $ws.Cells.Item(23, 10).Value2
$ws.Cells.Item(24, 11).Value = $ws.Cells.Item(23, 11).Value2
$ws.Cells.Item(24, 12).Value = $ws.Cells.Item(23, 12).Value2
$ws.Cells.Item(24, 13).Value = $ws.Cells.Item(23, 13).Value2

# Column F ("horacio_id"): split "q154,q155" into one id per row.
$ws.Cells.Item(23, 6).Value = "q154"
$ws.Cells.Item(24, 6).Value = "q155"

# Column H ("topic"): both rows share the same new topic label.
$ws.Cells.Item(23, 8).Value = "Information and Advice"
$ws.Cells.Item(24, 8).Value = "Information and Advice"

# Column G ("id"): each row gets its own unique id.
$ws.Cells.Item(23, 7).Value = "AJD_information"
$ws.Cells.Item(24, 7).Value = "AJD_inst_advice"

# Restore the view: scroll the frozen data pane back near the top and leave
# the selection on the newly added cell.
$ws.Activate()
$ws.Range("A18").Select()
$ws.Range("G25").Select()
